# Add data for 2022-07-30
#
# The workbook tracks Chicago "violent crime, year-to-date" counts, one
# worksheet per neighborhood (plus "Citywide Totals" and "By Neighborhood"
# roll-up sheets). Each sheet has crime-category rows (Aggravated Assault,
# Aggravated Battery, Criminal Sexual Assault, Homicide, Robbery, Total)
# down the side and years (2015-2022, some years omitted on sheets with no
# incidents that year) across the top. Adding one more day of data bumps
# the YTD count for that calendar day in every affected year/category/
# neighborhood combination by one.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 68
$ws.Range("C3").Value = 53
$ws.Range("E3").Value = 87
$ws.Range("I3").Value = 116
$ws.Range("C6").Value = 284
$ws.Range("D6").Value = 255
$ws.Range("G6").Value = 296
$ws.Range("H6").Value = 256
$ws.Range("I6").Value = 318
$ws.Range("C7").Value = 382
$ws.Range("D7").Value = 398
$ws.Range("E7").Value = 388
$ws.Range("G7").Value = 429
$ws.Range("H7").Value = 393
$ws.Range("I7").Value = 519

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("G19").Value = 14
$ws.Range("H32").Value = 30
$ws.Range("I47").Value = 15
$ws.Range("I53").Value = 84
$ws.Range("I65").Value = 15
$ws.Range("D74").Value = 8
$ws.Range("H74").Value = 10
$ws.Range("C76").Value = 10
$ws.Range("E77").Value = 17
$ws.Range("C80").Value = 5
$ws.Range("C88").Value = 6
$ws.Range("C96").Value = 8
$ws.Range("C98").Value = 382
$ws.Range("D98").Value = 398
$ws.Range("E98").Value = 388
$ws.Range("G98").Value = 429
$ws.Range("H98").Value = 393
$ws.Range("I98").Value = 519

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("C5").Value = 7
$ws.Range("C6").Value = 10

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("E3").Value = 4
$ws.Range("E7").Value = 17

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("H6").Value = 22
$ws.Range("H7").Value = 30

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("C3").Value = 2
$ws.Range("C5").Value = 6

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("C4").Value = 4
$ws.Range("C5").Value = 5

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("C5").Value = 7
$ws.Range("C6").Value = 8

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 9
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 3
$ws.Range("I6").Value = 15

$ws = $wb.Worksheets.Item("River North")
$ws.Range("D4").Value = 6
$ws.Range("H4").Value = 8
$ws.Range("D5").Value = 8
$ws.Range("H5").Value = 10

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("G5").Value = 8
$ws.Range("G6").Value = 14

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I5").Value = 13
$ws.Range("I6").Value = 15
